$wb = $excel.ActiveWorkbook

# ---------- Sheet 1: LP1912 ----------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = 'Última actualización: 09:42:42'
$ws1.Range("A3").Value = 'Total filas: 97'

$ws1.Range("A15").Value = '06:58:58'
$ws1.Range("B15").Value = '06:58'
$ws1.Range("C15").Value = '225_GOMEZ'
$ws1.Range("D15").Value = 0
$ws1.Range("E15").Value = 'LP1912'

$ws1.Range("A16").Value = '06:58:58'
$ws1.Range("B16").Value = '06:58'
$ws1.Range("C16").Value = '215A_EL PATO'
$ws1.Range("D16").Value = 0
$ws1.Range("E16").Value = 'LP1912'

$ws1.Range("A35").Value = '07:51:40'
$ws1.Range("B35").Value = '08:00'
$ws1.Range("C35").Value = '17_ROMERO'
$ws1.Range("D35").Value = 9
$ws1.Range("E35").Value = 'LP1912'

$ws1.Range("A36").Value = '07:26:49'
$ws1.Range("B36").Value = '08:00'
$ws1.Range("C36").Value = '16_SANTA ANA'
$ws1.Range("D36").Value = 34
$ws1.Range("E36").Value = 'LP1912'

$ws1.Range("A55").Value = '08:14:55'
$ws1.Range("B55").Value = '08:53'
$ws1.Range("C55").Value = '215B_EL PATO'
$ws1.Range("D55").Value = 39
$ws1.Range("E55").Value = 'LP1912'

$ws1.Range("A56").Value = '08:49:06'
$ws1.Range("B56").Value = '08:54'
$ws1.Range("C56").Value = '23_HERNANDEZ'
$ws1.Range("D56").Value = 4
$ws1.Range("E56").Value = 'LP1912'

$ws1.Range("A64").Value = '08:49:06'
$ws1.Range("B64").Value = '09:17'
$ws1.Range("C64").Value = '27_EL RETIRO'
$ws1.Range("D64").Value = 28
$ws1.Range("E64").Value = 'LP1912'

$ws1.Range("A65").Value = '08:57:42'
$ws1.Range("B65").Value = '09:17'
$ws1.Range("C65").Value = '14_ABASTO'
$ws1.Range("D65").Value = 20
$ws1.Range("E65").Value = 'LP1912'

$ws1.Range("A71").Value = '08:49:06'
$ws1.Range("B71").Value = '09:31'
$ws1.Range("C71").Value = '23_HERNANDEZ'
$ws1.Range("D71").Value = 42
$ws1.Range("E71").Value = 'LP1912'

$ws1.Range("A72").Value = '08:14:55'
$ws1.Range("B72").Value = '09:31'
$ws1.Range("C72").Value = '16_SANTA ANA'
$ws1.Range("D72").Value = 77
$ws1.Range("E72").Value = 'LP1912'

$ws1.Range("A78").Value = '09:42:42'
$ws1.Range("B78").Value = '09:42'
$ws1.Range("C78").Value = '11_ETCHEVERRY'
$ws1.Range("D78").Value = 0
$ws1.Range("E78").Value = 'LP1912'

$ws1.Range("A79").Value = '09:42:42'
$ws1.Range("B79").Value = '09:43'
$ws1.Range("C79").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D79").Value = 1
$ws1.Range("E79").Value = 'LP1912'

$ws1.Range("A80").Value = '09:42:42'
$ws1.Range("B80").Value = '09:53'
$ws1.Range("C80").Value = '10_OLMOS'
$ws1.Range("D80").Value = 11
$ws1.Range("E80").Value = 'LP1912'

$ws1.Range("A81").Value = '09:42:42'
$ws1.Range("B81").Value = '09:55'
$ws1.Range("C81").Value = '16_SANTA ANA'
$ws1.Range("D81").Value = 13
$ws1.Range("E81").Value = 'LP1912'

$ws1.Range("A82").Value = '09:42:42'
$ws1.Range("B82").Value = '09:58'
$ws1.Range("C82").Value = '215C_EL PATO'
$ws1.Range("D82").Value = 16
$ws1.Range("E82").Value = 'LP1912'

$ws1.Range("A83").Value = '08:49:06'
$ws1.Range("B83").Value = '09:59'
$ws1.Range("C83").Value = '215C_EL PATO'
$ws1.Range("D83").Value = 70
$ws1.Range("E83").Value = 'LP1912'

$ws1.Range("A84").Value = '08:57:42'
$ws1.Range("B84").Value = '10:05'
$ws1.Range("C84").Value = '14_ABASTO'
$ws1.Range("D84").Value = 68
$ws1.Range("E84").Value = 'LP1912'

$ws1.Range("A85").Value = '08:49:06'
$ws1.Range("B85").Value = '10:06'
$ws1.Range("C85").Value = '14_ABASTO'
$ws1.Range("D85").Value = 77
$ws1.Range("E85").Value = 'LP1912'

$ws1.Range("A86").Value = '09:42:42'
$ws1.Range("B86").Value = '10:13'
$ws1.Range("C86").Value = '17X38_ROMERO'
$ws1.Range("D86").Value = 31
$ws1.Range("E86").Value = 'LP1912'

$ws1.Range("A87").Value = '09:42:42'
$ws1.Range("B87").Value = '10:21'
$ws1.Range("C87").Value = '23_HERNANDEZ'
$ws1.Range("D87").Value = 39
$ws1.Range("E87").Value = 'LP1912'

$ws1.Range("A88").Value = '08:57:42'
$ws1.Range("B88").Value = '10:24'
$ws1.Range("C88").Value = '23_HERNANDEZ'
$ws1.Range("D88").Value = 87
$ws1.Range("E88").Value = 'LP1912'

$ws1.Range("A89").Value = '09:42:42'
$ws1.Range("B89").Value = '10:25'
$ws1.Range("C89").Value = '16_SANTA ANA'
$ws1.Range("D89").Value = 43
$ws1.Range("E89").Value = 'LP1912'

$ws1.Range("A90").Value = '09:42:42'
$ws1.Range("B90").Value = '10:29'
$ws1.Range("C90").Value = '15_ABASTO'
$ws1.Range("D90").Value = 47
$ws1.Range("E90").Value = 'LP1912'

$ws1.Range("A91").Value = '09:42:42'
$ws1.Range("B91").Value = '10:29'
$ws1.Range("C91").Value = '14_ABASTO'
$ws1.Range("D91").Value = 47
$ws1.Range("E91").Value = 'LP1912'

$ws1.Range("A92").Value = '09:42:42'
$ws1.Range("B92").Value = '10:44'
$ws1.Range("C92").Value = '11X44_ETCHEVERRY'
$ws1.Range("D92").Value = 62
$ws1.Range("E92").Value = 'LP1912'

$ws1.Range("A93").Value = '09:42:42'
$ws1.Range("B93").Value = '10:46'
$ws1.Range("C93").Value = '15_P INDUSTRIAL'
$ws1.Range("D93").Value = 64
$ws1.Range("E93").Value = 'LP1912'

$ws1.Range("A94").Value = '09:42:42'
$ws1.Range("B94").Value = '10:53'
$ws1.Range("C94").Value = '27_EL RETIRO'
$ws1.Range("D94").Value = 71
$ws1.Range("E94").Value = 'LP1912'

$ws1.Range("A95").Value = '09:42:42'
$ws1.Range("B95").Value = '10:59'
$ws1.Range("C95").Value = '10_OLMOS'
$ws1.Range("D95").Value = 77
$ws1.Range("E95").Value = 'LP1912'

$ws1.Range("A96").Value = '09:42:42'
$ws1.Range("B96").Value = '11:01'
$ws1.Range("C96").Value = '81_EL PELIGRO'
$ws1.Range("D96").Value = 79
$ws1.Range("E96").Value = 'LP1912'

$ws1.Range("A97").Value = '09:42:42'
$ws1.Range("B97").Value = '11:06'
$ws1.Range("C97").Value = '23_HERNANDEZ'
$ws1.Range("D97").Value = 84
$ws1.Range("E97").Value = 'LP1912'

$ws1.Range("A98").Value = '09:42:42'
$ws1.Range("B98").Value = '11:10'
$ws1.Range("C98").Value = '16_P MOR-SANTA ANA'
$ws1.Range("D98").Value = 88
$ws1.Range("E98").Value = 'LP1912'

$ws1.Range("A99").Value = '09:42:42'
$ws1.Range("B99").Value = '11:14'
$ws1.Range("C99").Value = '14_ABASTO'
$ws1.Range("D99").Value = 92
$ws1.Range("E99").Value = 'LP1912'

$ws1.Range("A100").Value = '09:42:42'
$ws1.Range("B100").Value = '11:15'
$ws1.Range("C100").Value = '15X38_ABASTO'
$ws1.Range("D100").Value = 93
$ws1.Range("E100").Value = 'LP1912'

$ws1.Range("A101").Value = '09:42:42'
$ws1.Range("B101").Value = '11:29'
$ws1.Range("C101").Value = '10_OLMOS'
$ws1.Range("D101").Value = 107
$ws1.Range("E101").Value = 'LP1912'

$ws1.Range("A102").Value = '09:42:42'
$ws1.Range("B102").Value = '11:30'
$ws1.Range("C102").Value = '215C_EL PATO'
$ws1.Range("D102").Value = 108
$ws1.Range("E102").Value = 'LP1912'

# ---------- Sheet 2: LP1912-215 ----------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = 'Última actualización: 09:42:42'
$ws2.Range("A3").Value = 'Total filas: 17'

$ws2.Range("A20").Value = '09:42:42'
$ws2.Range("B20").Value = '09:58'
$ws2.Range("C20").Value = '215C_EL PATO'
$ws2.Range("D20").Value = 16
$ws2.Range("E20").Value = 'LP1912'

$ws2.Range("A22").Value = '09:42:42'
$ws2.Range("B22").Value = '11:30'
$ws2.Range("C22").Value = '215C_EL PATO'
$ws2.Range("D22").Value = 108
$ws2.Range("E22").Value = 'LP1912'

# ---------- Sheet 3: 6203-6173 ----------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = 'Última actualización: 09:42:42'
$ws3.Range("A3").Value = 'Total filas: 17'

$ws3.Range("A16").Value = '09:42:42'
$ws3.Range("B16").Value = '10:12'
$ws3.Range("C16").Value = '215C_LA PLATA'
$ws3.Range("D16").Value = 30
$ws3.Range("E16").Value = 'L6203'

$ws3.Range("A18").Value = '09:42:42'
$ws3.Range("B18").Value = '10:29'
$ws3.Range("C18").Value = '215B_LP-P MOR-1 Y 57'
$ws3.Range("D18").Value = 47
$ws3.Range("E18").Value = 'L6173'

$ws3.Range("A19").Value = '08:49:06'
$ws3.Range("B19").Value = '10:30'
$ws3.Range("C19").Value = '215B_LP-P MOR-1 Y 57'
$ws3.Range("D19").Value = 101
$ws3.Range("E19").Value = 'L6173'

$ws3.Range("A20").Value = '09:42:42'
$ws3.Range("B20").Value = '10:30'
$ws3.Range("C20").Value = '215A_LA PLATA'
$ws3.Range("D20").Value = 48
$ws3.Range("E20").Value = 'L6173'

$ws3.Range("A22").Value = '09:42:42'
$ws3.Range("B22").Value = '11:25'
$ws3.Range("C22").Value = '215C_LA PLATA'
$ws3.Range("D22").Value = 103
$ws3.Range("E22").Value = 'L6203'

